$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: 2021-01-18, "9AM-12PM", "Core Java", "Done Assignment 1 In Java"
$ws.Range("B15").Copy() | Out-Null
$ws.Range("B16").PasteSpecial(-4122) | Out-Null
$ws.Range("B16").Value = 44214
$ws.Range("C16").Value = "9AM-12PM"
$ws.Range("D16").Value = "Core Java"
$ws.Range("E16").Value = "Done Assignment 1 In Java"

# Row 17: 2021-01-18, "1PM-10PM", "Core Java", "Building OOPS Application in JAVA"
$ws.Range("B15").Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4122) | Out-Null
$ws.Range("B17").Value = 44214
$ws.Range("D17").Value = "Core Java"
$ws.Range("E17").Value = "Building OOPS Application in JAVA"
$ws.Range("C17").Value = "1PM-10PM"

$excel.CutCopyMode = 0

$ws.Range("C17").Select() | Out-Null
